$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 383
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 383
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 383
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -733
$ws.Range("H28").Value = 1010.5
$ws.Range("I28").Value = 355.42856
$ws.Range("K28").Value = 355.42856
$ws.Range("M28").Value = 129.57144
$ws.Range("H33").Value = 178.61539
$ws.Range("I33").Value = 188.75
$ws.Range("K33").Value = 188.75
$ws.Range("M33").Value = 40.25
$ws.Range("H129").Value = 1746.7142
$ws.Range("J129").Value = 2670.5
$ws.Range("L129").Value = 8011.5
$ws.Range("N129").Value = -18011.5
$ws.Range("H132").Value = 5680.5864
$ws.Range("I132").Value = 5971.7407
$ws.Range("K132").Value = 17915.2221
$ws.Range("M132").Value = -15385.2221
$ws.Range("H137").Value = 1730266.1
$ws.Range("I137").Value = 2175856.2
$ws.Range("K137").Value = 6527568.600000001
$ws.Range("M137").Value = -6525018.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1279.5652
$ws.Range("I2").Value = 1198.7142
$ws.Range("K2").Value = 1198.7142
$ws.Range("M2").Value = -1085.7142
$ws.Range("H32").Value = 2732734.8
$ws.Range("I32").Value = 1236524.2
$ws.Range("K32").Value = 1236524.2
$ws.Range("M32").Value = -1236237.2
$ws.Range("H61").Value = 2685.0286
$ws.Range("I61").Value = 1821.5238
$ws.Range("J61").Value = 3980.2856
$ws.Range("K61").Value = 1821.5238
$ws.Range("L61").Value = 3980.2856
$ws.Range("M61").Value = -1609.5238
$ws.Range("N61").Value = -4404.2856
$ws.Range("H63").Value = 1496.6666
$ws.Range("I63").Value = 1496.6666
$ws.Range("K63").Value = 1496.6666
$ws.Range("M63").Value = -810.6666
$ws.Range("H66").Value = 1496.6666
$ws.Range("I66").Value = 1496.6666
$ws.Range("K66").Value = 7483.333000000001
$ws.Range("M66").Value = -4051.333000000001
$ws.Range("H116").Value = 1279.5652
$ws.Range("I116").Value = 1198.7142
$ws.Range("K116").Value = 1198.7142
$ws.Range("M116").Value = 1095.2858
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = ""
$ws.Range("H125").Value = 114907
$ws.Range("J125").Value = 114907
$ws.Range("L125").Value = 114907
$ws.Range("N125").Value = -124747
$ws.Range("H132").Value = 3071.3076
$ws.Range("I132").Value = 2792.7
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8378.099999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5848.099999999999
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 2685.0286
$ws.Range("I136").Value = 1821.5238
$ws.Range("J136").Value = 3980.2856
$ws.Range("K136").Value = 5464.5714
$ws.Range("L136").Value = 11940.8568
$ws.Range("M136").Value = -2914.5714
$ws.Range("N136").Value = -17040.8568
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 180000
$ws.Range("J141").Value = 180000
$ws.Range("L141").Value = 180000
$ws.Range("N141").Value = -190360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1279.5652
$ws.Range("I3").Value = 1198.7142
$ws.Range("K3").Value = 1198.7142
$ws.Range("M3").Value = -1084.7142
$ws.Range("H94").Value = 333350750
$ws.Range("J94").Value = 2500
$ws.Range("L94").Value = 2500
$ws.Range("N94").Value = -3402
$ws.Range("H132").Value = 104756
$ws.Range("J132").Value = 104756
$ws.Range("L132").Value = 104756
$ws.Range("N132").Value = -114876
$ws.Range("H133").Value = 107979.5
$ws.Range("J133").Value = 107979.5
$ws.Range("L133").Value = 107979.5
$ws.Range("N133").Value = -118099.5
$ws.Range("H134").Value = 2737.6924
$ws.Range("I134").Value = 2198.75
$ws.Range("K134").Value = 6596.25
$ws.Range("M134").Value = -4061.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3053118
$ws.Range("I31").Value = 3284.7222
$ws.Range("K31").Value = 3284.7222
$ws.Range("M31").Value = -2989.7222
$ws.Range("H34").Value = 3053118
$ws.Range("I34").Value = 3284.7222
$ws.Range("K34").Value = 3284.7222
$ws.Range("M34").Value = -3082.7222
$ws.Range("H86").Value = 4647.1665
$ws.Range("I86").Value = 3519.8
$ws.Range("J86").Value = 5452.4287
$ws.Range("K86").Value = 3519.8
$ws.Range("L86").Value = 5452.4287
$ws.Range("M86").Value = -2396.8
$ws.Range("N86").Value = -7698.4287
$ws.Range("H89").Value = 4647.1665
$ws.Range("I89").Value = 3519.8
$ws.Range("J89").Value = 5452.4287
$ws.Range("K89").Value = 17599
$ws.Range("L89").Value = 27262.1435
$ws.Range("M89").Value = -11983
$ws.Range("N89").Value = -38494.14350000001
$ws.Range("H105").Value = 2057.6667
$ws.Range("I105").Value = 1645.5714
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1645.5714
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = 101.4286
$ws.Range("N105").Value = -6994
$ws.Range("H132").Value = 3843.795
$ws.Range("I132").Value = 3479.6538
$ws.Range("K132").Value = 10438.9614
$ws.Range("M132").Value = -7908.9614
$ws.Range("H134").Value = 3685.75
$ws.Range("I134").Value = 4018.875
$ws.Range("J134").Value = 2686.375
$ws.Range("K134").Value = 12056.625
$ws.Range("L134").Value = 8059.125
$ws.Range("M134").Value = -9521.625
$ws.Range("N134").Value = -13129.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 23888.777
$ws.Range("J74").Value = 27495.5
$ws.Range("L74").Value = 82486.5
$ws.Range("N74").Value = -84608.5
$ws.Range("H77").Value = 23888.777
$ws.Range("J77").Value = 27495.5
$ws.Range("L77").Value = 247459.5
$ws.Range("N77").Value = -258067.5
$ws.Range("H107").Value = 437.41666
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("H114").Value = 1795.2307
$ws.Range("I114").Value = 1382.5555
$ws.Range("K114").Value = 4147.666499999999
$ws.Range("M114").Value = -893.6664999999994
$ws.Range("H140").Value = 2845.2856
$ws.Range("I140").Value = 2640.6875
$ws.Range("K140").Value = 7922.0625
$ws.Range("M140").Value = -2742.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2366.0264
$ws.Range("I132").Value = 2021.75
$ws.Range("J132").Value = 2616.4092
$ws.Range("K132").Value = 6065.25
$ws.Range("L132").Value = 7849.2276
$ws.Range("M132").Value = -3535.25
$ws.Range("N132").Value = -12909.2276
$ws.Range("H139").Value = 110999
$ws.Range("J139").Value = 110999
$ws.Range("L139").Value = 110999
$ws.Range("N139").Value = -121279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4886.3335
$ws.Range("J7").Value = 6124.5
$ws.Range("L7").Value = 6124.5
$ws.Range("N7").Value = -6348.5
$ws.Range("H14").Value = 10004
$ws.Range("I14").Value = 10004
$ws.Range("K14").Value = 10004
$ws.Range("M14").Value = -9832
$ws.Range("H40").Value = 20825.363
$ws.Range("J40").Value = 2569.5
$ws.Range("L40").Value = 2569.5
$ws.Range("N40").Value = -2841.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H126").Value = 4886.3335
$ws.Range("J126").Value = 6124.5
$ws.Range("L126").Value = 18373.5
$ws.Range("N126").Value = -23313.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4999
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""
$ws.Range("H86").Value = 32473
$ws.Range("J86").Value = 32473
$ws.Range("L86").Value = 32473
$ws.Range("N86").Value = -34719
$ws.Range("H89").Value = 32473
$ws.Range("J89").Value = 32473
$ws.Range("L89").Value = 162365
$ws.Range("N89").Value = -173597
$ws.Range("H100").Value = 83334060
$ws.Range("I100").Value = 907.7143
$ws.Range("K100").Value = 1815.4286
$ws.Range("M100").Value = -1274.4286
$ws.Range("H132").Value = 1426.4857
$ws.Range("I132").Value = 1364.5186
$ws.Range("K132").Value = 4093.5558
$ws.Range("M132").Value = -1563.5558
$ws.Range("H136").Value = 4030.375
$ws.Range("I136").Value = 2132.2334
$ws.Range("K136").Value = 6396.7002
$ws.Range("M136").Value = -3846.7002
